# Update column F (dSF) values per repull of data / mean calculation fix.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -2
    4  = 1
    5  = -1
    7  = 1
    14 = -2
    15 = 1
    25 = 4
    29 = 3
    33 = -1
    34 = -2
    44 = 1
    47 = 1
    51 = -5
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
